$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2 through 135 all change from 45178 (2023-09-09)
# to 45179 (2023-09-10). Update the serial date value directly so the
# existing date formatting (style) on the cells is preserved.
for ($r = 2; $r -le 135; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
